$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 6.584661960601807
$ws.Range("B1").Value = 5.37336254119873
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 3.843790769577026
$ws.Range("E1").Value = 1.904935359954834
